$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newAuthors = @(
  "[Elizabeth J.%Anderson%xref no email%1,           Kacey C.%Ernst%xref no email%1,           Francisco Fernando%Martins%xref no email%1,           Cicera da Silva%Martins%xref no email%1,           Mary P.%Koss%xref no email%1]",
  "[Ana Rosa%Linde-Arias%linde14@yahoo.com%1,                   Maria%Roura%NULL%1,                   Eduardo%Siqueira%NULL%1]",
  "[Maria%Belizan%mbelizan@iecs.org.ar%1,                   Edna%Maradiaga%edjamar3006@yahoo.com%1,                   Javier%Roberti%jroberti@iecs.org.ar%1,                   Maricela%Casco-Aguilar%marykasco@yahoo.com%1,                   Alison F.%Ortez%alison_fabiola@yahoo.es%1,                   Juan C.%Avila-Flores%javilaflores3@gmail.com%1,                   Gloria%González%marilyntoin@yahoo.com%1,                   Carolina%Bustillo%mcbu1502@yahoo.com%0,                   Alejandra%Calderón%lilianalecalderon@gmail.com%1,                   Harry%Bock%hbockme@hotmail.com%1,                   María L.%Cafferata%NULL%1,                   Adriano B.%Tavares%adriano.b.tavares@gmail.com%1,                   Jackeline%Alger%jackelinealger@gmail.com%1,                   Moazzam%Ali%alimoa@who.int%1]",
  "[ Catherine A.%Boyd%null%1,         Julie A.%Gazmararian%null%1,         Winifred Wilkins%Thompson%null%1]",
  "[Benner%P.%coreGivesNoEmail%1,                  Bohme%G.%coreGivesNoEmail%1,                  Chiang%H-H%coreGivesNoEmail%4,                  Chiang%H-H%coreGivesNoEmail%0,                  Chiang%H-H%coreGivesNoEmail%0,                  Chiang%H-H%coreGivesNoEmail%0,                  Donohoe%J.%coreGivesNoEmail%1,                  Hsien-Hsien%Chiang%coreGivesNoEmail%1,                  I-Ling%Sue%coreGivesNoEmail%1,                  International%Council of Nurses%coreGivesNoEmail%1,                  Mei-Bih%Chen%coreGivesNoEmail%1,                  Morgan%DL%coreGivesNoEmail%1,                  Packer%MJ%coreGivesNoEmail%1]",
  "[Chee%YO%coreGivesNoEmail%1,                  Dodgson%JE%coreGivesNoEmail%1,                  Tarrant%M%coreGivesNoEmail%1,                  Watkins%A%coreGivesNoEmail%1]",
  "[ Elin%Erland%null%1,                    Bente%Dahl%null%1,                  Elin%Erland%null%1,                  Bente%Dahl%null%1]",
  "[Davis%Mark%coreGivesNoEmail%1,                  Flowers%Paul%coreGivesNoEmail%1,                  Lohm%Davina%coreGivesNoEmail%1,                  Stephenson%Niamh%coreGivesNoEmail%1,                  Waller%Emily%coreGivesNoEmail%1]",
  "[Hector M.%Gomez%NULL%1,                   Carlos%Mejia Arbelaez%NULL%1,                   Jovana A.%Ocampo Cañas%ja.ocampo@uniandes.edu.co%1]",
  "[ Susan%Jones%null%1,                    Betty%Sam%null%1,                    Florence%Bull%null%1,                    Steven Bagie%Pieh%null%1,                    Jaki%Lambert%null%1,                    Florence%Mgawadere%null%1,                    Somasundari%Gopalakrishnan%null%1,                    Charles A.%Ameh%null%1,                    Nynke%van den Broek%null%1,                  Susan%Jones%null%1,                  Betty%Sam%null%1,                  Florence%Bull%null%1,                  Steven Bagie%Pieh%null%1,                  Jaki%Lambert%null%1,                  Florence%Mgawadere%null%1,                  Somasundari%Gopalakrishnan%null%1,                  Charles A.%Ameh%null%1,                  Nynke%van den Broek%null%1]",
  "[Caroline S.E.%Homer%NULL%1,                   Miranda%Davies-Tuck%NULL%1,                   Hannah G.%Dahlen%NULL%1,                   Vanessa L.%Scarf%NULL%1]",
  "[Ana Rosa%Linde%xref no email%1,           Carlos Eduardo%Siqueira%xref no email%1]",
  "[ Molly M.%Lynch%null%1,         Elizabeth W.%Mitchell%null%1,         Jennifer L.%Williams%null%1,         Kelly%Brumbaugh%null%1,         Michelle%Jones-Bell%null%1,         Debra E.%Pinkney%null%1,         Christine M.%Layton%null%1,         Patricia W.%Mersereau%null%1,         Juliette S.%Kendrick%null%1,         Paula Eguino%Medina%null%1,         Lucia Rojas%Smith%null%1]",
  "[Niuniu%Sun%NULL%1,                   Luoqun%Wei%NULL%1,                   Suling%Shi%NULL%1,                   Dandan%Jiao%NULL%1,                   Runluo%Song%NULL%1,                   Lili%Ma%NULL%1,                   Hongwei%Wang%NULL%1,                   Chao%Wang%NULL%1,                   Zhaoguo%Wang%NULL%1,                   Yanli%You%NULL%1,                   Shuhua%Liu%NULL%1,                   Hongyun%Wang%NULL%1]",
  "[Bushra%Sabri%xref no email%1,           Maria%Hartley%xref no email%1,           Jyoti%Saha%xref no email%1,           Sarah%Murray%xref no email%1,           Nancy%Glass%xref no email%1,           Jacquelyn C.%Campbell%xref no email%1]",
  "[Veronika%Tirado%Veronika.Tirado@ki.se%1,                   Santiago A.%Morales Mesa%NULL%1,                   John%Kinsman%NULL%1,                   Anna Mia%Ekström%NULL%1,                   Berta N.%Restrepo Jaramillo%NULL%1]",
  "[Simon N%Williams%NULL%1,                   Christopher J%Armitage%NULL%2,                   Christopher J%Armitage%NULL%0,                   Tova%Tampe%NULL%1,                   Kimberly%Dienes%NULL%1]"
)

for ($i = 0; $i -lt $newAuthors.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 5).Value = $newAuthors[$i]
}
